# Registration page and Registration page test
#
# The "users" sheet's USER_ID column (column A, rows 2-51) is refreshed with
# a new batch of generated usernames: same name part, numeric suffix doubled
# (e.g. Dayton1 -> Dayton2, Oscar2 -> Oscar4, Keshawn3 -> Keshawn6, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

$newUserIds = @(
    "Dayton2",
    "Oscar4",
    "Keshawn6",
    "Madden8",
    "Matteo10",
    "Tiffany12",
    "Odin14",
    "Dylan16",
    "Jaylee18",
    "Julie20",
    "Khalil22",
    "Jerry24",
    "Cecelia26",
    "Dax28",
    "Linda30",
    "Yusuf32",
    "Jenna34",
    "Ruth36",
    "Payten38",
    "Meredith40",
    "Trent42",
    "Camila44",
    "Kyson46",
    "Yadira48",
    "Lucille50",
    "Sanaa52",
    "Aileen54",
    "Leonard56",
    "April58",
    "Paulina60",
    "Nehemiah62",
    "Heath64",
    "Mikaela66",
    "Damari68",
    "Haylie70",
    "Raul72",
    "Kailey74",
    "Naima76",
    "Dakota78",
    "Shelby80",
    "Carter82",
    "Addyson84",
    "Mike86",
    "Irvin88",
    "Rocco90",
    "Julia92",
    "Karley94",
    "Dustin96",
    "Melanie98",
    "Lorelai100"
)

for ($i = 0; $i -lt $newUserIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newUserIds[$i]
}
